# Updates cryptos list - price (D) and volume% (E) columns, and two row re-orderings (B/C/D/E for rows 41/43 and 47/48)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "68.337.41"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "3.978.38"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'485.08"
$ws.Range("E5").Value = "  +8.24%  "
$ws.Range("D6").Value = "'149.81"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.737"
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  +8.56%  "
$ws.Range("D11").Value = "'0.0000371"
$ws.Range("E11").Value = "  +14.54%  "
$ws.Range("D12").Value = "'43.84"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "4.593.15"
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("D14").Value = "'10.53"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "'14.93"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "3.979.31"
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").Value = "'20.03"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "68.248.65"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'439.65"
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("D22").Value = "'3.41"
$ws.Range("E22").Value = "  +3.83%  "
$ws.Range("D23").Value = "'14.46"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").Value = "'88.55"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "'3.67"
$ws.Range("E25").Value = "  +5.73%  "
$ws.Range("D26").Value = "'39.00"
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("D27").Value = "'10.15"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "'732.18"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "'13.34"
$ws.Range("E30").Value = "  -3.73%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("D32").Value = "'2.84"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("D33").Value = "0.0₃0888"
$ws.Range("E33").Value = "  +28.93%  "
$ws.Range("D34").Value = "'42.12"
$ws.Range("E34").Value = "  -4.39%  "
$ws.Range("D35").Value = "'60.65"
$ws.Range("E35").Value = "  +6.54%  "
$ws.Range("D36").Value = "'0.152"
$ws.Range("E36").Value = "  -4.95%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'5.39"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("D39").Value = "'0.0475"
$ws.Range("E39").Value = "  -1.39%  "
$ws.Range("D40").Value = "'3.06"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'2.27"
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'2.88"
$ws.Range("E43").Value = "  +8.25%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'0.336"
$ws.Range("E45").Value = "  -0.63%  "
$ws.Range("D46").Value = "'3.45"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.27"
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "'2.54"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").Value = "'149.14"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "'2.89"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").Value = "'25.41"
$ws.Range("E51").Value = "  +0.37%  "
